# ------------------------------------------------------------------
# 1) Refresh the F-column ("time_taken") timestamps on the "data" sheet.
#    These reflect a re-run of the panelapp data pull.
# ------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$dataWs = $wb.Worksheets.Item("data")

$newTimes = @(
  "2021-10-05 14:35:29.416511",
  "2021-10-05 14:35:29.416519",
  "2021-10-05 14:35:29.416522",
  "2021-10-05 14:35:29.416525",
  "2021-10-05 14:35:29.416527",
  "2021-10-05 14:35:29.416530",
  "2021-10-05 14:35:29.416532",
  "2021-10-05 14:35:29.416535",
  "2021-10-05 14:35:29.416538",
  "2021-10-05 14:35:29.416540",
  "2021-10-05 14:35:29.416543",
  "2021-10-05 14:35:29.416545",
  "2021-10-05 14:35:29.416548",
  "2021-10-05 14:35:29.416550",
  "2021-10-05 14:35:29.416553",
  "2021-10-05 14:35:29.416555",
  "2021-10-05 14:35:29.416558",
  "2021-10-05 14:35:29.416560",
  "2021-10-05 14:35:29.416563",
  "2021-10-05 14:35:29.416565",
  "2021-10-05 14:35:29.416568",
  "2021-10-05 14:35:29.416570",
  "2021-10-05 14:35:29.416573",
  "2021-10-05 14:35:29.416575",
  "2021-10-05 14:35:29.416578",
  "2021-10-05 14:35:29.416580",
  "2021-10-05 14:35:29.416583",
  "2021-10-05 14:35:29.416585",
  "2021-10-05 14:35:29.416588",
  "2021-10-05 14:35:29.416590",
  "2021-10-05 14:35:29.416593",
  "2021-10-05 14:35:29.416595",
  "2021-10-05 14:35:29.416598",
  "2021-10-05 14:35:29.416600",
  "2021-10-05 14:35:29.416603",
  "2021-10-05 14:35:29.416605",
  "2021-10-05 14:35:29.416608",
  "2021-10-05 14:35:29.416610",
  "2021-10-05 14:35:29.416613",
  "2021-10-05 14:35:29.416615",
  "2021-10-05 14:35:29.416618",
  "2021-10-05 14:35:29.416620",
  "2021-10-05 14:35:29.416623",
  "2021-10-05 14:35:29.416625",
  "2021-10-05 14:35:29.416628",
  "2021-10-05 14:35:29.416630",
  "2021-10-05 14:35:29.416633",
  "2021-10-05 14:35:29.416635",
  "2021-10-05 14:35:29.416638",
  "2021-10-05 14:35:29.416640",
  "2021-10-05 14:35:29.416642",
  "2021-10-05 14:35:29.416645",
  "2021-10-05 14:35:29.416648",
  "2021-10-05 14:35:29.416650",
  "2021-10-05 14:35:29.416653",
  "2021-10-05 14:35:29.416655",
  "2021-10-05 14:35:29.416658",
  "2021-10-05 14:35:29.416660",
  "2021-10-05 14:35:29.416663",
  "2021-10-05 14:35:29.416665",
  "2021-10-05 14:35:29.416668",
  "2021-10-05 14:35:29.416670",
  "2021-10-05 14:35:29.416672",
  "2021-10-05 14:35:29.416675",
  "2021-10-05 14:35:29.416678",
  "2021-10-05 14:35:29.416681",
  "2021-10-05 14:35:29.416684",
  "2021-10-05 14:35:29.416686",
  "2021-10-05 14:35:29.416689",
  "2021-10-05 14:35:29.416691",
  "2021-10-05 14:35:29.416693",
  "2021-10-05 14:35:29.416696",
  "2021-10-05 14:35:29.416698",
  "2021-10-05 14:35:29.416701",
  "2021-10-05 14:35:29.416703",
  "2021-10-05 14:35:29.416706",
  "2021-10-05 14:35:29.416710",
  "2021-10-05 14:35:29.416713",
  "2021-10-05 14:35:29.416716",
  "2021-10-05 14:35:29.416718",
  "2021-10-05 14:35:29.416720",
  "2021-10-05 14:35:29.416723",
  "2021-10-05 14:35:29.416725",
  "2021-10-05 14:35:29.416728",
  "2021-10-05 14:35:29.416730",
  "2021-10-05 14:35:29.416733",
  "2021-10-05 14:35:29.416735",
  "2021-10-05 14:35:29.416738",
  "2021-10-05 14:35:29.416740",
  "2021-10-05 14:35:29.416742",
  "2021-10-05 14:35:29.416745",
  "2021-10-05 14:35:29.416747",
  "2021-10-05 14:35:29.416751",
  "2021-10-05 14:35:29.416754",
  "2021-10-05 14:35:29.416756",
  "2021-10-05 14:35:29.416759",
  "2021-10-05 14:35:29.416761",
  "2021-10-05 14:35:29.416764",
  "2021-10-05 14:35:29.416766",
  "2021-10-05 14:35:29.416768",
  "2021-10-05 14:35:29.416771",
  "2021-10-05 14:35:29.416773",
  "2021-10-05 14:35:29.416776",
  "2021-10-05 14:35:29.416778"
)

for ($i = 0; $i -lt $newTimes.Length; $i++) {
  $row = $i + 2
  $dataWs.Cells.Item($row, 6).Value = $newTimes[$i]
}

# ------------------------------------------------------------------
# 2) Add a new "metadata" sheet (placed right after "data") that
#    records the panel query this pull corresponds to.
# ------------------------------------------------------------------
$metaWs = $wb.Worksheets.Add($null, $wb.Worksheets.Item("data"))
$metaWs.Name = "metadata"

$metaWs.Range("B1").Value = "data_name"
$metaWs.Range("C1").Value = "data_id"
$metaWs.Range("D1").Value = "data_version"
$metaWs.Range("E1").Value = "data_version_created"
$metaWs.Range("F1").Value = "panel_query_time"
$metaWs.Range("G1").Value = "panel_get_request"

$metaWs.Range("B1:G1").Font.Bold = $true
$metaWs.Range("B1:G1").Borders.LineStyle = 1
$metaWs.Range("B1:G1").HorizontalAlignment = -4108
$metaWs.Range("B1:G1").VerticalAlignment = -4160

$metaWs.Range("A2").Value = 0
$metaWs.Range("A2").Font.Bold = $true
$metaWs.Range("A2").Borders.LineStyle = 1
$metaWs.Range("A2").HorizontalAlignment = -4108
$metaWs.Range("A2").VerticalAlignment = -4160

$metaWs.Range("B2").Value = "Red cell disorders"
$metaWs.Range("C2").Value = 3366
$metaWs.Range("D2").NumberFormat = "@"
$metaWs.Range("D2").Value = "1.2"
$metaWs.Range("E2").Value = "2021-09-19T08:04:48.820183Z"
$metaWs.Range("F2").Value = "2021-10-05 14:35:29.413038"
$metaWs.Range("G2").Value = "https://panelapp.agha.umccr.org/api/v1/panels/3366/?format=json"

# Keep "data" as the active sheet/tab, matching the original workbook view.
$dataWs.Select()
